# The uploaded workbook replaces the first time-slot pair ("2:55-3:0" /
# "3:0-3:5") in column C with a new pair ("8:20-8:25" / "8:25-8:30"), while
# every other row (10:55-11:0, 11:0-11:5, 18:55-19:0, 19:0-19:5 and the two
# hashed "Посыл / Заповедь" values in column B) stays exactly as it was.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "8:20-8:25"
$ws.Range("C3").Value = "8:25-8:30"

# The saved file also shows the active selection moved from B11 to B12.
[void]$ws.Range("B12").Select()
